$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the "year" header to "season_ending_year"
$ws.Range("B1").Value = "season_ending_year"

# 2. Populate the previously-empty birth_year column (E) for all data rows
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 5).Value = 1937
}

# 3. Add a new "calendar_year" column (AM) mirroring the header style used
#    by the other header cells, and fill in the per-row calendar year values
$ws.Range("AM1").Value = "calendar_year"
$ws.Range("AL1").Copy()
$ws.Range("AM1").PasteSpecial(-4122)  # xlPasteFormats

$calendarYears = @{
    2  = 1973
    3  = 1972
    4  = 1971
    5  = 1970
    6  = 1969
    7  = 1968
    8  = 1967
    9  = 1966
    10 = 1965
    11 = 1965
    12 = 1965
    13 = 1964
    14 = 1963
    15 = 1962
    16 = 1961
    17 = 1960
}

foreach ($r in $calendarYears.Keys) {
    $ws.Cells.Item($r, 39).Value = $calendarYears[$r]
}
